# Add new row under range 0-150 in Card1 by admin
#
# A new row is inserted at row 3 of the "Card1" sheet (pushing the
# existing rows 3-12 down to rows 4-13), recording a new service event
# for the same Min_Tones/Max_Tones bucket (0-150) as row 2.
# Additionally, the previously-blank "Crrection" (P) column on the other
# data rows is populated with the placeholder text "nan".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# --- Insert a new row at position 3, shifting existing rows 3-12 down to 4-13 ---
$ws.Rows.Item(3).Insert()

# --- Populate the new row 3 (keep everything as literal text, like the rest of the sheet) ---
$newRow = $ws.Range("A3:P3")
$newRow.NumberFormat = "@"

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "0"
$ws.Range("C3").Value = "150"
$ws.Range("D3").Value = ".1"
# E3:L3 and P3 are intentionally left blank
$ws.Range("M3").Value = "تم تركيب وعيار ماكينه"
$ws.Range("N3").Value = "تم التشغيل "
$ws.Range("O3").Value = "م.الشناوي"

# --- Backfill the "Crrection" (P) column placeholder on the other rows ---
# (P3 is the brand-new row and stays blank; P2 and P4:P13 are the
#  pre-existing rows, set separately since COM multi-area ranges here
#  only apply .Value to the first area)
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "nan"
$ws.Range("P4:P13").NumberFormat = "@"
$ws.Range("P4:P13").Value = "nan"

Write-Output "done"
